$d = $word.ActiveDocument

# The "Requisitos" section ends with the paragraph:
#   "LOQ4083: Fenômenos de Transporte I (Requisito fraco)"
# followed by an empty paragraph, a page-break paragraph, and the site's
# copyright/footer paragraph. The site rebuild removed that trailing
# footer block (the blank paragraph, the page-break paragraph that
# precedes it, and the copyright paragraph itself), leaving the
# "LOQ4083..." paragraph followed directly by the final blank paragraph
# and page-break paragraph that close out the document.
#
# Locate the requirement paragraph by its text, then delete the three
# paragraphs that immediately follow it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "LOQ4083:*Requisito fraco*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Delete the three paragraphs right after the requirement paragraph:
    # 1) blank paragraph
    # 2) page-break blank paragraph
    # 3) "© 2020 . Contact: ..." copyright paragraph
    $d.Paragraphs.Item($target + 1).Range.Delete()
    $d.Paragraphs.Item($target + 1).Range.Delete()
    $d.Paragraphs.Item($target + 1).Range.Delete()
}
